$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("DBD")
$ws.Range("D19").Value = "TIMESTAMP"
$ws.Range("D21").Value = "TIMESTAMP"
